# Add season-record columns (Wins, Losses, Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1): Wins / Losses / Ties in columns AD:AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) by
# copying the style from the last existing header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill the season record for every player row (2-43): every row gets the
# same team record (92 wins, 70 losses, 0 ties).
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 92
    $ws.Cells.Item($r, 31).Value = 70
    $ws.Cells.Item($r, 32).Value = 0
}
